# Auto-generated Excel COM-interop edit script
# Applies updated market-board derived values (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the Leve profit tables across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Values come from a scheduled data-refresh run; only numeric <v> cell contents change,
# no structural/formatting changes. One special case: ARM row 102 col N is cleared
# entirely (its prior value collapses into the recomputed M102 figure).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 69445350
$ws.Range("I92").Value = 3968983.2
$ws.Range("J92").Value = 222223550
$ws.Range("K92").Value = 3968983.2
$ws.Range("L92").Value = 222223550
$ws.Range("M92").Value = -3967735.2
$ws.Range("N92").Value = -222226046

$ws.Range("H106").Value = 74077784
$ws.Range("I106").Value = 22226276
$ws.Range("J106").Value = 333335330
$ws.Range("K106").Value = 22226276
$ws.Range("L106").Value = 333335330
$ws.Range("M106").Value = -22225645
$ws.Range("N106").Value = -333336592

$ws.Range("H113").Value = 17859444
$ws.Range("I113").Value = 2684.8333
$ws.Range("K113").Value = 2684.8333
$ws.Range("M113").Value = 569.1667000000002

$ws.Range("H127").Value = 1446.3529
$ws.Range("I127").Value = 796.75
$ws.Range("J127").Value = 1646.2307
$ws.Range("K127").Value = 2390.25
$ws.Range("L127").Value = 4938.6921
$ws.Range("M127").Value = 2569.75
$ws.Range("N127").Value = -14858.6921

$ws.Range("H132").Value = 1080.7333
$ws.Range("I132").Value = 660.4838999999999
$ws.Range("K132").Value = 1981.4517
$ws.Range("M132").Value = 548.5483000000002

$ws.Range("H137").Value = 1267.8657
$ws.Range("I137").Value = 1058.16
$ws.Range("J137").Value = 1884.6471
$ws.Range("K137").Value = 3174.48
$ws.Range("L137").Value = 5653.9413
$ws.Range("M137").Value = -624.4800000000005
$ws.Range("N137").Value = -10753.9413

$ws.Range("H138").Value = 2151.4683
$ws.Range("I138").Value = 827.32654
$ws.Range("J138").Value = 4314.2334
$ws.Range("K138").Value = 2481.97962
$ws.Range("L138").Value = 12942.7002
$ws.Range("M138").Value = 2658.02038
$ws.Range("N138").Value = -23222.7002


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2255.55
$ws.Range("I32").Value = 2155.9404
$ws.Range("J32").Value = 2778.5
$ws.Range("K32").Value = 2155.9404
$ws.Range("L32").Value = 2778.5
$ws.Range("M32").Value = -1868.9404
$ws.Range("N32").Value = -3352.5

$ws.Range("H61").Value = 1294.902
$ws.Range("I61").Value = 1153.4736
$ws.Range("J61").Value = 1708.3077
$ws.Range("K61").Value = 1153.4736
$ws.Range("L61").Value = 1708.3077
$ws.Range("M61").Value = -941.4736
$ws.Range("N61").Value = -2132.3077

$ws.Range("H74").Value = 935.4423
$ws.Range("I74").Value = 863.8857400000001
$ws.Range("K74").Value = 863.8857400000001
$ws.Range("M74").Value = 10.11425999999994

$ws.Range("H77").Value = 935.4423
$ws.Range("I77").Value = 863.8857400000001
$ws.Range("K77").Value = 4319.4287
$ws.Range("M77").Value = 48.57129999999961

$ws.Range("H102").Value = 4631595
$ws.Range("I102").Value = 4631595
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4631595
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4629973
$ws.Range("N102").Value = ""

$ws.Range("H124").Value = 28141.166
$ws.Range("J124").Value = 28141.166
$ws.Range("L124").Value = 28141.166
$ws.Range("N124").Value = -37961.166

$ws.Range("H132").Value = 1615121.4
$ws.Range("I132").Value = 1514.4694
$ws.Range("J132").Value = 7697178
$ws.Range("K132").Value = 4543.4082
$ws.Range("L132").Value = 23091534
$ws.Range("M132").Value = -2013.4082
$ws.Range("N132").Value = -23096594

$ws.Range("H136").Value = 1294.902
$ws.Range("I136").Value = 1153.4736
$ws.Range("J136").Value = 1708.3077
$ws.Range("K136").Value = 3460.4208
$ws.Range("L136").Value = 5124.9231
$ws.Range("M136").Value = -910.4207999999999
$ws.Range("N136").Value = -10224.9231


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1664.55
$ws.Range("I86").Value = 1606.5714
$ws.Range("J86").Value = 1799.8334
$ws.Range("K86").Value = 1606.5714
$ws.Range("L86").Value = 1799.8334
$ws.Range("M86").Value = -483.5714
$ws.Range("N86").Value = -4045.8334

$ws.Range("H89").Value = 1664.55
$ws.Range("I89").Value = 1606.5714
$ws.Range("J89").Value = 1799.8334
$ws.Range("K89").Value = 8032.857
$ws.Range("L89").Value = 8999.166999999999
$ws.Range("M89").Value = -2416.857
$ws.Range("N89").Value = -20231.167

$ws.Range("H134").Value = 1422.0725
$ws.Range("I134").Value = 1030.5834
$ws.Range("J134").Value = 2316.9048
$ws.Range("K134").Value = 3091.7502
$ws.Range("L134").Value = 6950.714399999999
$ws.Range("M134").Value = -556.7501999999999
$ws.Range("N134").Value = -12020.7144

$ws.Range("H135").Value = 48887.184
$ws.Range("J135").Value = 48887.184
$ws.Range("L135").Value = 48887.184
$ws.Range("N135").Value = -59027.184


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 16780
$ws.Range("J20").Value = 16780
$ws.Range("L20").Value = 16780
$ws.Range("N20").Value = -17252

$ws.Range("H30").Value = 16780
$ws.Range("J30").Value = 16780
$ws.Range("L30").Value = 16780
$ws.Range("N30").Value = -16962

$ws.Range("H31").Value = 4457.846
$ws.Range("I31").Value = 1078.8182
$ws.Range("K31").Value = 1078.8182
$ws.Range("M31").Value = -783.8181999999999

$ws.Range("H34").Value = 4457.846
$ws.Range("I34").Value = 1078.8182
$ws.Range("K34").Value = 1078.8182
$ws.Range("M34").Value = -876.8181999999999

$ws.Range("H128").Value = 16780
$ws.Range("J128").Value = 16780
$ws.Range("L128").Value = 16780
$ws.Range("N128").Value = -26740

$ws.Range("H130").Value = 52520
$ws.Range("J130").Value = 52520
$ws.Range("L130").Value = 52520
$ws.Range("N130").Value = -62560

$ws.Range("H132").Value = 1678.2909
$ws.Range("I132").Value = 1311.175
$ws.Range("J132").Value = 2657.2666
$ws.Range("K132").Value = 3933.525
$ws.Range("L132").Value = 7971.7998
$ws.Range("M132").Value = -1403.525
$ws.Range("N132").Value = -13031.7998

$ws.Range("H134").Value = 1947.2858
$ws.Range("I134").Value = 2139.3635
$ws.Range("J134").Value = 1502.4736
$ws.Range("K134").Value = 6418.0905
$ws.Range("L134").Value = 4507.4208
$ws.Range("M134").Value = -3883.0905
$ws.Range("N134").Value = -9577.4208


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3226619.5
$ws.Range("I131").Value = 5556144.5
$ws.Range("J131").Value = 1123.1538
$ws.Range("K131").Value = 16668433.5
$ws.Range("L131").Value = 3369.4614
$ws.Range("M131").Value = -16663393.5
$ws.Range("N131").Value = -13449.4614


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 46298224
$ws.Range("I122").Value = 53242572
$ws.Range("J122").Value = 2566
$ws.Range("K122").Value = 159727716
$ws.Range("L122").Value = 7698
$ws.Range("M122").Value = -159725266
$ws.Range("N122").Value = -12598

$ws.Range("H132").Value = 1708.3281
$ws.Range("I132").Value = 1492.7727
$ws.Range("J132").Value = 2182.55
$ws.Range("K132").Value = 4478.3181
$ws.Range("L132").Value = 6547.650000000001
$ws.Range("M132").Value = -1948.3181
$ws.Range("N132").Value = -11607.65


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7321400
$ws.Range("I132").Value = 9058204
$ws.Range("J132").Value = 2014.2858
$ws.Range("K132").Value = 27174612
$ws.Range("L132").Value = 6042.857400000001
$ws.Range("M132").Value = -27172082
$ws.Range("N132").Value = -11102.8574

$ws.Range("H136").Value = 3223.349
$ws.Range("I136").Value = 1497.3518
$ws.Range("J136").Value = 13579.333
$ws.Range("K136").Value = 4492.055399999999
$ws.Range("L136").Value = 40737.999
$ws.Range("M136").Value = -1942.055399999999
$ws.Range("N136").Value = -45837.999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 789.2308
$ws.Range("I126").Value = 503.45456
$ws.Range("K126").Value = 1510.36368
$ws.Range("M126").Value = 959.6363200000001

$ws.Range("H132").Value = 23076
$ws.Range("I132").Value = 26992.5
$ws.Range("J132").Value = 1815
$ws.Range("K132").Value = 80977.5
$ws.Range("L132").Value = 5445
$ws.Range("M132").Value = -78447.5
$ws.Range("N132").Value = -10505

$ws.Range("H136").Value = 8774623
$ws.Range("I136").Value = 2758.4102
$ws.Range("J136").Value = 27780330
$ws.Range("K136").Value = 8275.230599999999
$ws.Range("L136").Value = 83340990
$ws.Range("M136").Value = -5725.230599999999
$ws.Range("N136").Value = -83346090


Write-Output "Applied scheduled market-data refresh to all sheets."
